$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell E8 text from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the selection change recorded in the workbook (user selected E8)
$ws.Range("E8").Select()
